$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1432.3684
$ws.Range("I98").Value = 666.3333
$ws.Range("J98").Value = 4305
$ws.Range("K98").Value = 666.3333
$ws.Range("L98").Value = 4305
$ws.Range("M98").Value = 831.6667
$ws.Range("N98").Value = -7301

$ws.Range("H107").Value = 37247.85
$ws.Range("I107").Value = 38642
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 38642
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -36722
$ws.Range("N107").Value = -4840

$ws.Range("H122").Value = 1432.3684
$ws.Range("I122").Value = 666.3333
$ws.Range("J122").Value = 4305
$ws.Range("K122").Value = 1998.9999
$ws.Range("L122").Value = 12915
$ws.Range("M122").Value = 451.0001
$ws.Range("N122").Value = -17815

$ws.Range("H137").Value = 1298.6129
$ws.Range("I137").Value = 1145.9565
$ws.Range("J137").Value = 1737.5
$ws.Range("K137").Value = 3437.8695
$ws.Range("L137").Value = 5212.5
$ws.Range("M137").Value = -887.8694999999998
$ws.Range("N137").Value = -10312.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 30000
$ws.Range("J42").Value = 30000
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -30972

$ws.Range("H52").Value = 12871.429
$ws.Range("J52").Value = 12871.429
$ws.Range("L52").Value = 12871.429
$ws.Range("N52").Value = -13507.429

$ws.Range("H102").Value = 2377.5
$ws.Range("I102").Value = 2377.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2377.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -755.5
$ws.Range("N102").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 535.9167
$ws.Range("I64").Value = 435.5
$ws.Range("J64").Value = 586.125
$ws.Range("K64").Value = 435.5
$ws.Range("L64").Value = 586.125
$ws.Range("M64").Value = -210.5
$ws.Range("N64").Value = -1036.125

$ws.Range("H67").Value = 535.9167
$ws.Range("I67").Value = 435.5
$ws.Range("J67").Value = 586.125
$ws.Range("K67").Value = 435.5
$ws.Range("L67").Value = 586.125
$ws.Range("M67").Value = 344.5
$ws.Range("N67").Value = -2146.125

$ws.Range("H81").Value = 16570
$ws.Range("J81").Value = 16570
$ws.Range("L81").Value = 16570
$ws.Range("N81").Value = -18692

$ws.Range("H84").Value = 16570
$ws.Range("J84").Value = 16570
$ws.Range("L84").Value = 49710
$ws.Range("N84").Value = -60318

$ws.Range("H105").Value = 2948
$ws.Range("I105").Value = 1875.4166
$ws.Range("K105").Value = 1875.4166
$ws.Range("M105").Value = -128.4166

$ws.Range("H107").Value = 3709.9285
$ws.Range("I107").Value = 2661.6667
$ws.Range("J107").Value = 9999.5
$ws.Range("K107").Value = 2661.6667
$ws.Range("L107").Value = 9999.5
$ws.Range("M107").Value = -741.6667000000002
$ws.Range("N107").Value = -13839.5

$ws.Range("H135").Value = 61000
$ws.Range("J135").Value = 61000
$ws.Range("L135").Value = 61000
$ws.Range("N135").Value = -71140

$ws.Range("H140").Value = 55000
$ws.Range("J140").Value = 55000
$ws.Range("L140").Value = 55000
$ws.Range("N140").Value = -65360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 150.33333
$ws.Range("I7").Value = 56.6
$ws.Range("J7").Value = 267.5
$ws.Range("K7").Value = 56.6
$ws.Range("L7").Value = 267.5
$ws.Range("M7").Value = 56.4
$ws.Range("N7").Value = -493.5

$ws.Range("H22").Value = 340
$ws.Range("I22").Value = 340
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 340
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 10
$ws.Range("N22").Value = $null

$ws.Range("H31").Value = 2214.2986
$ws.Range("I31").Value = 1499.079
$ws.Range("J31").Value = 2911.1794
$ws.Range("K31").Value = 1499.079
$ws.Range("L31").Value = 2911.1794
$ws.Range("M31").Value = -1204.079
$ws.Range("N31").Value = -3501.1794

$ws.Range("H34").Value = 2214.2986
$ws.Range("I34").Value = 1499.079
$ws.Range("J34").Value = 2911.1794
$ws.Range("K34").Value = 1499.079
$ws.Range("L34").Value = 2911.1794
$ws.Range("M34").Value = -1297.079
$ws.Range("N34").Value = -3315.1794

$ws.Range("H55").Value = 5718.25
$ws.Range("I55").Value = 5436.5
$ws.Range("K55").Value = 5436.5
$ws.Range("M55").Value = -5121.5

$ws.Range("H107").Value = 949.9666999999999
$ws.Range("I107").Value = 911
$ws.Range("J107").Value = 1008.4167
$ws.Range("K107").Value = 911
$ws.Range("L107").Value = 1008.4167
$ws.Range("M107").Value = 1009
$ws.Range("N107").Value = -4848.4167

$ws.Range("H122").Value = 2584.2727
$ws.Range("I122").Value = 1747.8
$ws.Range("J122").Value = 3281.3333
$ws.Range("K122").Value = 5243.4
$ws.Range("L122").Value = 9843.999899999999
$ws.Range("M122").Value = -2793.4
$ws.Range("N122").Value = -14743.9999

$ws.Range("H130").Value = 69980
$ws.Range("J130").Value = 69980
$ws.Range("L130").Value = 69980
$ws.Range("N130").Value = -80020

$ws.Range("H134").Value = 7387.7
$ws.Range("I134").Value = 9290.267
$ws.Range("J134").Value = 1680
$ws.Range("K134").Value = 27870.801
$ws.Range("L134").Value = 5040
$ws.Range("M134").Value = -25335.801
$ws.Range("N134").Value = -10110

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2501.7144
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2501.7144
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 7505.1432
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = -9877.143199999999

$ws.Range("H89").Value = 2501.7144
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2501.7144
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 22515.4296
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = -34371.4296

$ws.Range("H113").Value = 568.45
$ws.Range("I113").Value = 502.7143
$ws.Range("J113").Value = 603.8461
$ws.Range("K113").Value = 1508.1429
$ws.Range("L113").Value = 1811.5383
$ws.Range("M113").Value = 661.8571000000002
$ws.Range("N113").Value = -6151.5383

$ws.Range("H121").Value = 741350.1
$ws.Range("I121").Value = 247
$ws.Range("J121").Value = 1667729
$ws.Range("K121").Value = 741
$ws.Range("L121").Value = 5003187
$ws.Range("M121").Value = 569
$ws.Range("N121").Value = -5005807

$ws.Range("H140").Value = 26059.512
$ws.Range("I140").Value = 48400.906
$ws.Range("J140").Value = 2601.05
$ws.Range("K140").Value = 145202.718
$ws.Range("L140").Value = 7803.150000000001
$ws.Range("M140").Value = -140022.718
$ws.Range("N140").Value = -18163.15

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4138.5
$ws.Range("I122").Value = 3372.5715
$ws.Range("K122").Value = 10117.7145
$ws.Range("M122").Value = -7667.7145

$ws.Range("H140").Value = 56857.145
$ws.Range("J140").Value = 56857.145
$ws.Range("L140").Value = 56857.145
$ws.Range("N140").Value = -67217.14499999999

$ws.Range("H141").Value = 55685.8
$ws.Range("J141").Value = 55685.8
$ws.Range("L141").Value = 55685.8
$ws.Range("N141").Value = -66045.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = $null

$ws.Range("H40").Value = 2245.7144
$ws.Range("I40").Value = 1807.7273
$ws.Range("J40").Value = 3851.6667
$ws.Range("K40").Value = 1807.7273
$ws.Range("L40").Value = 3851.6667
$ws.Range("M40").Value = -1671.7273
$ws.Range("N40").Value = -4123.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3281.0476
$ws.Range("I132").Value = 3209.2812
$ws.Range("J132").Value = 3510.7
$ws.Range("K132").Value = 9627.8436
$ws.Range("L132").Value = 10532.1
$ws.Range("M132").Value = -7097.8436
$ws.Range("N132").Value = -15592.1

$ws.Range("H135").Value = 50713.125
$ws.Range("J135").Value = 50713.125
$ws.Range("L135").Value = 50713.125
$ws.Range("N135").Value = -60853.125

$ws.Range("H140").Value = 40164.5
$ws.Range("J140").Value = 40164.5
$ws.Range("L140").Value = 40164.5
$ws.Range("N140").Value = -50524.5

$ws.Range("H141").Value = 47607.5
$ws.Range("J141").Value = 47607.5
$ws.Range("L141").Value = 47607.5
$ws.Range("N141").Value = -57967.5
